$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsmatrix")

# --- Row 43: new entry "Systemarchitektur Anforderungen sammeln" ---
$ws.Range("A43").Value = 9
$ws.Range("B43").Value = "Konzeptuelles Design"
$ws.Range("C43").Value = "[TASK]"
$ws.Range("D43").Value = "Systemarchitektur"
$ws.Range("E43").Value = "Systemarchitektur Anforderungen sammeln"
$ws.Range("F43").Value = 44320
$ws.Range("G43").Value = 44338
$ws.Range("J43").Value = 0.375
$ws.Range("K43").Value = 0.54166666666666663

# --- Row 44: new entry "Systemarchitektur modellieren" ---
$ws.Range("A44").Value = 9
$ws.Range("B44").Value = "Konzeptuelles Design"
$ws.Range("C44").Value = "[TASK]"
$ws.Range("D44").Value = "Systemarchitektur"
$ws.Range("E44").Value = "Systemarchitektur modellieren"
$ws.Range("F44").Value = 44320
$ws.Range("G44").Value = 44338
$ws.Range("J44").Value = 0.58333333333333337
$ws.Range("K44").Value = 0.78125

# Extend the ROUNDUP "hours worked" formula down into the two new rows
# (I15:I42 already share formula si=0; this creates the si=1 group for I43:I44)
$ws.Range("I43:I44").Formula = "=ROUNDUP(((SUM(K43-J43)*24*60/60)/0.25),0)*0.25"

# --- Row 45: F45/G45 go back to being completely empty (no cell, no style) ---
$ws.Range("F45:G45").Clear()

# --- Data validation: row 45 no longer participates in the Prefix list rule ---
$ws.Range("C45").Validation.Delete()

# --- Sheet view / selection housekeeping ---
$win = $excel.ActiveWindow
$win.ScrollRow = 24
$win.ScrollColumn = 1
[void]$ws.Range("H53").Select()
